# Updates crypto Price (D) and Volume(1h) (E) columns on Sheet1 with the
# latest scraped values, as produced by the GitHub Actions symbol-list
# update job.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new text value. Values are written as text (not
# numbers) because the source sheet stores every Price/Volume cell as an
# inline string (e.g. "303.39", "-1.63%"), and the COM layer would
# otherwise auto-coerce numeric-/percent-looking literals into real
# numbers. Forcing the cell to a text number-format before the write keeps
# it a genuine string, matching the original sheet's cell typing.
$updates = [ordered]@{
    "D2"  = "303.39";     "E2"  = "-1.63%"
    "D3"  = "35.71";      "E3"  = "-1.37%"
    "D4"  = "5.039";      "E4"  = "-1.26%"
    "D5"  = "0.07944";    "E5"  = "-2.22%"
    "D6"  = "1.859";      "E6"  = "-4.41%"
                          "E7"  = "-0.60%"
                          "E8"  = "-1.25%"
    "D9"  = "0.1344";     "E9"  = "-3.39%"
    "D10" = "0.1879";     "E10" = "-2.42%"
    "D11" = "0.09042";    "E11" = "-2.38%"
    "D12" = "0.03442";    "E12" = "1.61%"
    "D13" = "0.09813";    "E13" = "-0.52%"
    "D14" = "0.001401";   "E14" = "-1.46%"
    "D15" = "0.006053";   "E15" = "5.38%"
    "D16" = "3.734";      "E16" = "3.09%"
    "D17" = "4.103";      "E17" = "-1.99%"
                          "E18" = "12.46%"
    "D19" = "0.3444";     "E19" = "0.17%"
    "D20" = "0.1332";     "E20" = "-1.27%"
    "D21" = "5.158";      "E21" = "5.27%"
    "D22" = "0.2397";     "E22" = "-4.09%"
    "D23" = "0.04388";    "E23" = "-2.72%"
    "D24" = "0.001232";   "E24" = "1.10%"
    "D25" = "0.004617";   "E25" = "-5.19%"
    "D26" = "0.0001293";  "E26" = "4.26%"
    "D27" = "0.0004436";  "E27" = "-0.12%"
    "D39" = "0.01931";    "E39" = "-4.50%"
    "D40" = "0.05247";    "E40" = "5.92%"
    "D41" = "0.007587";   "E41" = "-0.83%"
    "D42" = "0.01012";    "E42" = "-1.44%"
    "D43" = "0.1347";     "E43" = "-2.74%"
    "D44" = "0.002149";   "E44" = "2.28%"
    "D45" = "0.01012";    "E45" = "-10.68%"
    "D46" = "0.00006135"; "E46" = "-4.61%"
                          "E47" = "-0.15%"
    "D48" = "65.22";      "E48" = "0.85%"
    "D49" = "0.001657";   "E49" = "39.12%"
    "D50" = "0.00002098"; "E50" = "-0.15%"
    "D51" = "0.0001998";  "E51" = "-0.15%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
